# Insert a new data row for 2014-03-30 (20:30 -> 22:00, 90 min / 1.5 h)
# right where the previous blank spacer row (row 103) used to be. This
# pushes the old blank spacer row and the three summary rows down by one,
# and Excel keeps the SUM()/aggregate formulas pointing at the right ranges
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 103 (and everything below it) down by inserting a fresh row.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new time entry.
$ws.Range("A103").Value = 2014
$ws.Range("B103").Value = 3
$ws.Range("C103").Value = 30
$ws.Range("D103").Value = 0.85416666666666663
$ws.Range("E103").Value = 0.91666666666666663
$ws.Range("F103").Formula = "=(E103-D103)*24*60"
$ws.Range("G103").Formula = "=F103/60"

# Match the workbook's saved selection (F103) after the edit.
$ws.Range("F103").Select()
